{"js": "// Update the two-digit-divided-by-one-digit division problems in the\n// table to a new set of operands, per the commit's regenerated output.\n// Each \"old\" expression occurs exactly once in the document, so a\n// plain text search-and-replace on the run text is sufficient and\n// preserves all existing run/paragraph formatting.\nconst replacements = [\n  [\"36\u00f79=\", \"34\u00f79=\"],\n  [\"19\u00f76=\", \"53\u00f74=\"],\n  [\"51\u00f72=\", \"98\u00f75=\"],\n  [\"58\u00f79=\", \"89\u00f72=\"],\n  [\"67\u00f75=\", \"89\u00f73=\"],\n  [\"62\u00f77=\", \"18\u00f75=\"],\n  [\"46\u00f72=\", \"76\u00f76=\"],\n  [\"17\u00f74=\", \"30\u00f75=\"],\n  [\"12\u00f76=\", \"50\u00f79=\"],\n  [\"51\u00f79=\", \"71\u00f75=\"],\n  [\"25\u00f76=\", \"62\u00f75=\"],\n  [\"99\u00f74=\", \"83\u00f74=\"],\n  [\"68\u00f72=\", \"36\u00f77=\"],\n  [\"11\u00f76=\", \"70\u00f74=\"],\n  [\"48\u00f76=\", \"48\u00f72=\"],\n  [\"54\u00f78=\", \"77\u00f79=\"],\n  [\"59\u00f78=\", \"28\u00f74=\"],\n  [\"65\u00f79=\", \"59\u00f76=\"],\n  [\"30\u00f77=\", \"29\u00f78=\"],\n  [\"92\u00f76=\", \"89\u00f73=\"],\n  [\"17\u00f75=\", \"92\u00f72=\"],\n  [\"88\u00f79=\", \"73\u00f79=\"],\n  [\"70\u00f77=\", \"95\u00f74=\"],\n  [\"57\u00f75=\", \"99\u00f75=\"],\n  [\"87\u00f77=\", \"65\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-divided-by-one-digit division problems in the\n# table to a new set of operands, per the commit's regenerated output.\n# Each \"old\" expression occurs exactly once in the document, so a\n# plain Find/Replace on the whole document content is sufficient and\n# preserves all existing run/paragraph formatting.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"36\u00f79=\", \"34\u00f79=\"),\n  @(\"19\u00f76=\", \"53\u00f74=\"),\n  @(\"51\u00f72=\", \"98\u00f75=\"),\n  @(\"58\u00f79=\", \"89\u00f72=\"),\n  @(\"67\u00f75=\", \"89\u00f73=\"),\n  @(\"62\u00f77=\", \"18\u00f75=\"),\n  @(\"46\u00f72=\", \"76\u00f76=\"),\n  @(\"17\u00f74=\", \"30\u00f75=\"),\n  @(\"12\u00f76=\", \"50\u00f79=\"),\n  @(\"51\u00f79=\", \"71\u00f75=\"),\n  @(\"25\u00f76=\", \"62\u00f75=\"),\n  @(\"99\u00f74=\", \"83\u00f74=\"),\n  @(\"68\u00f72=\", \"36\u00f77=\"),\n  @(\"11\u00f76=\", \"70\u00f74=\"),\n  @(\"48\u00f76=\", \"48\u00f72=\"),\n  @(\"54\u00f78=\", \"77\u00f79=\"),\n  @(\"59\u00f78=\", \"28\u00f74=\"),\n  @(\"65\u00f79=\", \"59\u00f76=\"),\n  @(\"30\u00f77=\", \"29\u00f78=\"),\n  @(\"92\u00f76=\", \"89\u00f73=\"),\n  @(\"17\u00f75=\", \"92\u00f72=\"),\n  @(\"88\u00f79=\", \"73\u00f79=\"),\n  @(\"70\u00f77=\", \"95\u00f74=\"),\n  @(\"57\u00f75=\", \"99\u00f75=\"),\n  @(\"87\u00f77=\", \"65\u00f73=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
